# Add "2022-Q4" worksheet data, inserted right after the "总计" summary sheet.
# This mirrors a quarterly refresh: a brand-new "2022-Q4" tab is added, all
# other quarter tabs shift one position to the right (unchanged otherwise),
# and the "总计" (totals) sheet gains a new leading row summarizing the
# new quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q4" worksheet right after "总计".
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q4 = $wb.Worksheets.Add($null, $zongji)
$q4.Name = "2022-Q4"

# Header row (same layout/style as the other quarterly sheets).
$header = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($i = 0; $i -lt $header.Length; $i++) {
    $q4.Cells.Item(1, 2 + $i).Value2 = $header[$i]
}
$headerRange = $q4.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fund holdings data for 2022-Q4.
$rows = @(
    @("516150", "嘉实中证稀土产业ETF",         "20.83", "99.33", "3.51", "0.7311", 10),
    @("516780", "华泰柏瑞中证稀土产业ETF",       "7.87",  "98.83", "3.51", "0.2762", 10),
    @("159713", "富国中证稀土产业ETF",           "2.62",  "99.35", "3.50", "0.0917", 10),
    @("159715", "易方达中证稀土产业ETF",         "2.52",  "98.82", "3.47", "0.0874", 10),
    @("015061", "中信建投沪深300指数增强A",      "1.22",  "89.86", "1.78", "0.0217", 5),
    @("015062", "中信建投沪深300指数增强C",      "0.81",  "89.86", "1.78", "0.0144", 5)
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    $excelRow = 2 + $r

    $a = $q4.Cells.Item($excelRow, 1)
    $a.Value2 = $r
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Font.Bold = $true
    $a.Borders.LineStyle = 1

    $q4.Cells.Item($excelRow, 2).Value2 = "'" + $row[0]
    $q4.Cells.Item($excelRow, 3).Value2 = $row[1]
    $q4.Cells.Item($excelRow, 4).Value2 = "'" + $row[2]
    $q4.Cells.Item($excelRow, 5).Value2 = "'" + $row[3]
    $q4.Cells.Item($excelRow, 6).Value2 = "'" + $row[4]
    $q4.Cells.Item($excelRow, 7).Value2 = "'" + $row[5]
    $q4.Cells.Item($excelRow, 8).Value2 = $row[6]
}

# ---------------------------------------------------------------------
# 2) Update the "总计" (totals) sheet: insert the 2022-Q4 summary row at
#    the top of the data (row 2), pushing the existing rows down.
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")

$summary = @(
    @("2022-Q4", 6,  1.22),
    @("2022-Q3", 4,  1.25),
    @("2022-Q2", 4,  1.61),
    @("2022-Q1", 6,  1.75),
    @("2021-Q4", 4,  2.01),
    @("2021-Q2", 6,  2.8),
    @("2021-Q1", 15, 2.4),
    @("2020-Q4", 41, 14.03)
)

for ($r = 0; $r -lt $summary.Length; $r++) {
    $row = $summary[$r]
    $excelRow = 2 + $r

    $a = $zongji.Cells.Item($excelRow, 1)
    $a.Value2 = $r
    $a.HorizontalAlignment = -4108
    $a.VerticalAlignment = -4160
    $a.Font.Bold = $true
    $a.Borders.LineStyle = 1

    $zongji.Cells.Item($excelRow, 2).Value2 = $row[0]
    $zongji.Cells.Item($excelRow, 3).Value2 = $row[1]
    $zongji.Cells.Item($excelRow, 4).Value2 = $row[2]
}

$zongji.Select()
$zongji.Range("A1").Select()
